{"js": "// Applies the four changes from the diff:\n// 1. Remove the stray \"_GoBack\" bookmark from the title paragraph.\n// 2. \"...record the posterior mean estimate of...\" -> \"...record the estimate of...\"\n// 3. \"confidence interval for the posterior mean parameter.\" ->\n//    \"interval for the mean parameter.\" (keeping the two-run split, just\n//    shifted: run B becomes \"interval for the \", run C becomes \"mean parameter.\")\n// 4. \"Calculate the proportion of replicates for which the true mean was\n//    within the 80% credible interval.\" -> drop \"credible \" and re-insert the\n//    \"_GoBack\" bookmark right before the final \"interval.\"\n\nconst body = context.document.body;\n\n// --- 1. Remove the old _GoBack bookmark -----------------------------------\ncontext.document.deleteBookmark('_GoBack');\nawait context.sync();\n\n// --- 2. \"posterior mean estimate\" -> \"estimate\" ----------------------------\nconst estimateResults = body.search('record the posterior mean estimate of', { matchCase: true });\nestimateResults.load('items');\nawait context.sync();\nif (estimateResults.items.length > 0) {\n  estimateResults.items[0].insertText('record the estimate of', 'Replace');\n  await context.sync();\n}\n\n// --- 3. \"confidence interval for the posterior mean parameter.\" -----------\n// Locate the paragraph that still contains \"confidence\" so the two edits\n// below stay scoped to it (keeps the two runs involved unambiguous).\nconst paragraphs = body.paragraphs;\nparagraphs.load('items');\nawait context.sync();\nparagraphs.items.forEach((p) => p.load('text'));\nawait context.sync();\n\nlet confidenceParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf('confidence') !== -1) {\n    confidenceParagraph = p;\n    break;\n  }\n}\n\nif (confidenceParagraph) {\n  // 3a. Run B: \"confidence \" -> \"interval for the \"\n  const confRange = confidenceParagraph.getRange();\n  const confResults = confRange.search('confidence ', { matchCase: true });\n  confResults.load('items');\n  await context.sync();\n  if (confResults.items.length > 0) {\n    confResults.items[0].insertText('interval for the ', 'Replace');\n    await context.sync();\n  }\n\n  // 3b. Run C: \"interval for the posterior mean parameter.\" -> \"mean parameter.\"\n  const tailRange = confidenceParagraph.getRange();\n  const tailResults = tailRange.search('interval for the posterior mean parameter.', { matchCase: true });\n  tailResults.load('items');\n  await context.sync();\n  if (tailResults.items.length > 0) {\n    tailResults.items[0].insertText('mean parameter.', 'Replace');\n    await context.sync();\n  }\n}\n\n// --- 4. \"80% credible interval.\" -> \"80% \" + bookmark + \"interval.\" -------\nconst credibleResults = body.search('80% credible interval.', { matchCase: true });\ncredibleResults.load('items');\nawait context.sync();\nif (credibleResults.items.length > 0) {\n  credibleResults.items[0].insertText('80% interval.', 'Replace');\n  await context.sync();\n}\n\nconst intervalResults = body.search('interval.', { matchCase: true });\nintervalResults.load('items');\nawait context.sync();\nif (intervalResults.items.length > 0) {\n  const startRange = intervalResults.items[0].getRange('Start');\n  startRange.insertBookmark('_GoBack');\n  await context.sync();\n}\n", "ps1": "# Applies the four changes from the diff:\n# 1. Remove the stray \"_GoBack\" bookmark from the title paragraph.\n# 2. \"...record the posterior mean estimate of...\" -> \"...record the estimate of...\"\n# 3. \"confidence interval for the posterior mean parameter.\" ->\n#    \"interval for the mean parameter.\" (two runs keep existing boundary,\n#    just the text each side of it changes).\n# 4. \"Calculate the proportion of replicates for which the true mean was\n#    within the 80% credible interval.\" -> drop \"credible \" and re-insert the\n#    \"_GoBack\" bookmark right before the final \"interval.\"\n#\n# NOTE: this COM host coalesces adjacent runs that end up with identical\n# formatting as soon as any of their text is touched. Several of the runs\n# on either side of our edits share formatting, so before touching any\n# text we drop temporary \"fence\" bookmarks at every run boundary we need\n# to keep intact; they get removed again at the very end.\n\n$d = $word.ActiveDocument\n\n# --- 1. Remove the old _GoBack bookmark ------------------------------------\nif ($d.Bookmarks.Exists('_GoBack')) {\n    $d.Bookmarks('_GoBack').Delete()\n}\n\n# --- locate the two target paragraphs --------------------------------------\n$fitParagraph = $null\n$calcParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like '*confidence*') {\n        $fitParagraph = $p\n    }\n    if ($p.Range.Text -like '*credible interval*') {\n        $calcParagraph = $p\n    }\n}\n\n# --- fence every run boundary inside $fitParagraph we must not disturb -----\n$fenceNames = @('FENCE_A', 'FENCE_B', 'FENCE_C', 'FENCE_D', 'FENCE_E', 'FENCE_F')\n$fenceAnchors = @(\n    'You can use any',\n    '2-2 afternoon -- intro to linear models.R',\n    [char]8221 + '.  Also, please record the',\n    ' 80% ',\n    'confidence ',\n    'interval for the posterior mean parameter.'\n)\n\nfor ($i = 0; $i -lt $fenceAnchors.Length; $i++) {\n    $fr = $fitParagraph.Range\n    $ok = $fr.Find.Execute($fenceAnchors[$i], $true)\n    $fr.Collapse(1)\n    $d.Bookmarks.Add($fenceNames[$i], $fr)\n}\n\n# --- 2. \"posterior mean estimate\" -> \"estimate\" -----------------------------\n$r1 = $fitParagraph.Range\n$r1.Find.Execute('posterior mean estimate of', $true, $false, $false, $false, $false, $true, 1, $false, 'estimate of', 1)\n\n# --- 3a. run: \"confidence \" -> \"interval for the \" --------------------------\n$r2 = $fitParagraph.Range\n$r2.Find.Execute('confidence ', $true, $false, $false, $false, $false, $true, 1, $false, 'interval for the ', 1)\n\n# --- 3b. run: \"interval for the posterior mean parameter.\" -> \"mean parameter.\" ---\n$r3 = $fitParagraph.Range\n$r3.Find.Execute('interval for the posterior mean parameter.', $true, $false, $false, $false, $false, $true, 1, $false, 'mean parameter.', 1)\n\n# --- remove the temporary fence bookmarks -----------------------------------\nforeach ($name in $fenceNames) {\n    if ($d.Bookmarks.Exists($name)) {\n        $d.Bookmarks($name).Delete()\n    }\n}\n\n# --- 4. \"80% credible interval.\" -> \"80% \" + bookmark + \"interval.\" --------\n$r4 = $calcParagraph.Range\n$r4.Find.Execute('credible ', $true, $false, $false, $false, $false, $true, 1, $false, '', 1)\n\n$r5 = $calcParagraph.Range\n$r5.Find.Execute('interval.', $true)\n$r5.Collapse(1)\n$d.Bookmarks.Add('_GoBack', $r5)\n"}
